$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1) labels for columns AG:AS
$ws.Range("AG1").Value = "A"
$ws.Range("AH1").Value = "B"
$ws.Range("AI1").Value = "C"
$ws.Range("AJ1").Value = "D"
$ws.Range("AK1").Value = "E"
$ws.Range("AL1").Value = "F"
$ws.Range("AM1").Value = "G"
$ws.Range("AN1").Value = "H"
$ws.Range("AO1").Value = "Finlande,etc"
$ws.Range("AP1").Value = "1000+"
$ws.Range("AQ1").Value = "1000-"
$ws.Range("AR1").Value = "W"
$ws.Range("AS1").Value = "T"

# Row 2 data
$ws.Range("AG2").Value = 0.7
$ws.Range("AH2").Value = 0.7
$ws.Range("AI2").Value = 0.7
$ws.Range("AJ2").Value = 0.7
$ws.Range("AK2").Value = 1
$ws.Range("AL2").Value = 0.7
$ws.Range("AM2").Value = 0.7
$ws.Range("AN2").Value = 0
$ws.Range("AO2").Value = 0.7
$ws.Range("AP2").Value = 0.7
$ws.Range("AQ2").Value = 0.5
$ws.Range("AR2").Value = 0.6
$ws.Range("AS2").Value = 0.6

# Row 3 data
$ws.Range("AG3").Value = 0.5
$ws.Range("AH3").Value = 0.5
$ws.Range("AI3").Value = 0.7
$ws.Range("AJ3").Value = 0.7
$ws.Range("AK3").Value = 0.9
$ws.Range("AL3").Value = 0.7
$ws.Range("AM3").Value = 0.5
$ws.Range("AN3").Value = 0
$ws.Range("AO3").Value = 0.5
$ws.Range("AP3").Value = 0.5
$ws.Range("AQ3").Value = 0.2
$ws.Range("AR3").Value = 0.2
$ws.Range("AS3").Value = 0.5

# Row 4 data
$ws.Range("AG4").Value = 0.3
$ws.Range("AH4").Value = 0.3
$ws.Range("AI4").Value = 0.6
$ws.Range("AJ4").Value = 0.6
$ws.Range("AK4").Value = 0.8
$ws.Range("AL4").Value = 0.6
$ws.Range("AM4").Value = 0.3
$ws.Range("AN4").Value = 0
$ws.Range("AO4").Value = 0.2
$ws.Range("AP4").Value = 0.2
$ws.Range("AQ4").Value = 0
$ws.Range("AR4").Value = 0
$ws.Range("AS4").Value = 0

# Restore the view: zoom to the newly added columns and move the selection
$excel.ActiveWindow.Zoom = 130
$excel.ActiveWindow.ScrollColumn = 41
$excel.ActiveWindow.ScrollRow = 1
$ws.Range("AK8").Select()
